$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / info block text updates ---
$ws.Range("A1").Value = "Bieu mau Test"
$ws.Range("B2").Value = "Test"
$ws.Range("D2").Value = "Chỉ huyện"
$ws.Range("B3").Value = "Tháng"
$ws.Range("B4").Value = "Báo cáo năm 2020"

# --- Indicator table updates (rows 8-10) ---
$ws.Range("A8").Value = "Chi tiêu cha"
$ws.Range("C8").Value = "Bich"
$ws.Range("F8").Value = 1

$ws.Range("A9").Value = "Chi tieu con"
$ws.Range("C9").Value = "Hop"
$ws.Range("F9").Value = 2

$ws.Range("A10").Value = "Chi tieu test"
$ws.Range("C10").Value = "Hop"
$ws.Range("F10").Value = 3

# --- Row 11: fill in previously-empty indicator cells ---
$ws.Range("A11").Value = "Chi tieu test 5"
$ws.Range("C11").Value = "Hop"
$ws.Range("F11").Value = 4

# --- Add a new blank template row 12 (copy formatting from row 11) ---
$ws.Range("A11:C11").Copy($ws.Range("A12:C12"))
$ws.Range("A12").Value = ""
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""

# --- Restore selection to match the expanded data range ---
$ws.Range("A8:C12").Select() | Out-Null
